# "Generate Report for Handback"
#
# For each language sheet (zh-cn, de-de):
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - "Latest Target File" (F) / "Latest Handback File" (G) columns get filled in,
#     linking to the same source .md / handoff .xlf files as columns A / D
#   - "Latest Handback DateTime" (H) is stamped with the handback time

$wb = $excel.ActiveWorkbook

$mdFile  = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/a3bdd5a84487e6724dc3b7938fc67b89ed4db127/e2e/52ca86cf-fe3f-49d9-8e54-fafdcc507556.md"

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276 # RGB(0x64,0x95,0xED) == xl style "HyperLink" font color FF6495ED

# Hyperlinks.Add() already applies underline; just match the workbook's
# existing cornflower-blue hyperlink font color (instead of the default theme color).
function Apply-HyperlinkStyle($range) {
    $range.Font.Color = $hyperlinkColor
}

# ---------------- Overview sheet ----------------
# B2/C2/B3/C3 mirror the per-language Status text and need the same update.
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $statusNew
$ov.Range("C2").Value = $statusNew
$ov.Range("B3").Value = $statusNew
$ov.Range("C3").Value = $statusNew

# ---------------- zh-cn sheet ----------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

$zhXlf = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2bad171c62fe36dde12b5ebd632be1412304b5e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.zh-cn.xlf"

foreach ($row in 2,3) {
    $fCell = $zh.Range("F$row")
    $fCell.Value = $mdFile
    $zh.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdFile) | Out-Null
    Apply-HyperlinkStyle $fCell

    $gCell = $zh.Range("G$row")
    $gCell.Value = $zhXlf
    $zh.Hyperlinks.Add($gCell, $zhXlfUrl, "", "", $zhXlf) | Out-Null
    Apply-HyperlinkStyle $gCell
}

$zh.Range("H2").Value = "2016-03-22 13:13:24"
$zh.Range("H3").Value = "2016-03-22 13:13:24"

# ---------------- de-de sheet ----------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

$deXlf = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50589fd186930f33dbd543d920c7cb5c14a8e44f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.de-de.xlf"

foreach ($row in 2,3) {
    $fCell = $de.Range("F$row")
    $fCell.Value = $mdFile
    $de.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdFile) | Out-Null
    Apply-HyperlinkStyle $fCell

    $gCell = $de.Range("G$row")
    $gCell.Value = $deXlf
    $de.Hyperlinks.Add($gCell, $deXlfUrl, "", "", $deXlf) | Out-Null
    Apply-HyperlinkStyle $gCell
}

$de.Range("H2").Value = "2016-03-22 13:13:31"
$de.Range("H3").Value = "2016-03-22 13:13:31"
